# Auto-generated script to apply Halicarnassus_Profits market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 539.6667
$ws.Range("I53").Value = 539.6667
$ws.Range("K53").Value = 539.6667
$ws.Range("M53").Value = 97.33330000000001
$ws.Range("H86").Value = 6282.3335
$ws.Range("I86").Value = 6938.8
$ws.Range("K86").Value = 6938.8
$ws.Range("M86").Value = -5815.8
$ws.Range("H88").Value = 1799.5
$ws.Range("I88").Value = 1599.5
$ws.Range("K88").Value = 1599.5
$ws.Range("M88").Value = -1193.5
$ws.Range("H89").Value = 6282.3335
$ws.Range("I89").Value = 6938.8
$ws.Range("K89").Value = 34694
$ws.Range("M89").Value = -29078
$ws.Range("H91").Value = 1799.5
$ws.Range("I91").Value = 1599.5
$ws.Range("K91").Value = 1599.5
$ws.Range("M91").Value = -195.5
$ws.Range("H116").Value = 2653.1667
$ws.Range("I116").Value = 2546.625
$ws.Range("J116").Value = 2866.25
$ws.Range("K116").Value = 2546.625
$ws.Range("L116").Value = 2866.25
$ws.Range("M116").Value = 895.375
$ws.Range("N116").Value = -9750.25
$ws.Range("H132").Value = 13318.571
$ws.Range("I132").Value = 12205
$ws.Range("K132").Value = 36615
$ws.Range("M132").Value = -34085
$ws.Range("H135").Value = 1266
$ws.Range("I135").Value = 899.5
$ws.Range("J135").Value = 1999
$ws.Range("K135").Value = 8095.5
$ws.Range("L135").Value = 17991
$ws.Range("M135").Value = -5560.5
$ws.Range("N135").Value = -23061
$ws.Range("H138").Value = 2904.25
$ws.Range("I138").Value = 996
$ws.Range("K138").Value = 2988
$ws.Range("M138").Value = 2152
$ws.Range("H141").Value = 3683.3333
$ws.Range("J141").Value = 5550
$ws.Range("L141").Value = 16650
$ws.Range("N141").Value = -27010

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3507.862
$ws.Range("I61").Value = 2629.16
$ws.Range("K61").Value = 2629.16
$ws.Range("M61").Value = -2417.16
$ws.Range("H122").Value = 495.58334
$ws.Range("I122").Value = 495.58334
$ws.Range("K122").Value = 1486.75002
$ws.Range("M122").Value = 963.2499800000001
$ws.Range("H136").Value = 3507.862
$ws.Range("I136").Value = 2629.16
$ws.Range("K136").Value = 7887.48
$ws.Range("M136").Value = -5337.48

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1820
$ws.Range("I99").Value = 1421.5625
$ws.Range("K99").Value = 1421.5625
$ws.Range("M99").Value = 76.4375
$ws.Range("H134").Value = 3921.3333
$ws.Range("I134").Value = 3874.3333
$ws.Range("K134").Value = 11622.9999
$ws.Range("M134").Value = -9087.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4199.864
$ws.Range("I31").Value = 1967.5
$ws.Range("K31").Value = 1967.5
$ws.Range("M31").Value = -1672.5
$ws.Range("H34").Value = 4199.864
$ws.Range("I34").Value = 1967.5
$ws.Range("K34").Value = 1967.5
$ws.Range("M34").Value = -1765.5
$ws.Range("H58").Value = 5020.75
$ws.Range("I58").Value = 4663
$ws.Range("K58").Value = 4663
$ws.Range("M58").Value = -4460
$ws.Range("H62").Value = 1467.1428
$ws.Range("I62").Value = 1467.1428
$ws.Range("K62").Value = 1467.1428
$ws.Range("M62").Value = -843.1428000000001
$ws.Range("H65").Value = 1467.1428
$ws.Range("I65").Value = 1467.1428
$ws.Range("K65").Value = 7335.714
$ws.Range("M65").Value = -4215.714
$ws.Range("H98").Value = 62390
$ws.Range("J98").Value = 62390
$ws.Range("L98").Value = 62390
$ws.Range("N98").Value = -66882
$ws.Range("H100").Value = 45000
$ws.Range("J100").Value = 45000
$ws.Range("L100").Value = 45000
$ws.Range("N100").Value = -47164
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1996
$ws.Range("I132").Value = 1996
$ws.Range("K132").Value = 5988
$ws.Range("M132").Value = -3458
$ws.Range("H136").Value = 5020.75
$ws.Range("I136").Value = 4663
$ws.Range("K136").Value = 13989
$ws.Range("M136").Value = -11439

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 291
$ws.Range("I23").Value = 225
$ws.Range("J23").Value = 335
$ws.Range("K23").Value = 675
$ws.Range("L23").Value = 1005
$ws.Range("M23").Value = -440
$ws.Range("N23").Value = -1475

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2853.5
$ws.Range("I102").Value = 2710.4
$ws.Range("K102").Value = 2710.4
$ws.Range("M102").Value = -1088.4
$ws.Range("H132").Value = 4997.5
$ws.Range("I132").Value = 4997.5
$ws.Range("K132").Value = 14992.5
$ws.Range("M132").Value = -12462.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 877
$ws.Range("I22").Value = 666.3333
$ws.Range("K22").Value = 666.3333
$ws.Range("M22").Value = -371.3333
$ws.Range("H27").Value = 877
$ws.Range("I27").Value = 666.3333
$ws.Range("K27").Value = 666.3333
$ws.Range("M27").Value = -559.3333
$ws.Range("H40").Value = 2991.5833
$ws.Range("J40").Value = 1802
$ws.Range("L40").Value = 1802
$ws.Range("N40").Value = -2074
$ws.Range("H42").Value = 39999
$ws.Range("J42").Value = 39999
$ws.Range("L42").Value = 39999
$ws.Range("N42").Value = -41125
$ws.Range("H49").Value = 39999
$ws.Range("J49").Value = 39999
$ws.Range("L49").Value = 39999
$ws.Range("N49").Value = -40293
$ws.Range("H87").Value = 47950
$ws.Range("J87").Value = 47950
$ws.Range("L87").Value = 47950
$ws.Range("N87").Value = -50196
$ws.Range("H90").Value = 47950
$ws.Range("J90").Value = 47950
$ws.Range("L90").Value = 143850
$ws.Range("N90").Value = -155082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10000.25
$ws.Range("J62").Value = 12000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 10000.25
$ws.Range("J65").Value = 12000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240
$ws.Range("H107").Value = 2152.6365
$ws.Range("I107").Value = 2528.7778
$ws.Range("K107").Value = 7586.3334
$ws.Range("M107").Value = -5666.3334
